# Apply updated odds values to the FlashScore weekly games worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 2.4
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("AI2").Value = 13

# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 2.45
$ws.Range("J3").Value = 3.5
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 29
$ws.Range("AC3").Value = 8.5
$ws.Range("AD3").Value = 6
$ws.Range("AH3").Value = 7.5
$ws.Range("AK3").Value = 23
$ws.Range("AO3").Value = 17
$ws.Range("BB3").Value = 81

# Row 4
$ws.Range("G4").Value = 2.2
$ws.Range("I4").Value = 3
$ws.Range("L4").Value = 3.6
$ws.Range("AI4").Value = 15
$ws.Range("AJ4").Value = 11
$ws.Range("AL4").Value = 23
$ws.Range("AM4").Value = 29
$ws.Range("BB4").Value = 67

# Row 5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
